$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.073.09"
$ws.Range("E2").Value = "  -0.50%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.652.27"
$ws.Range("E3").Value = "  -0.46%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "217.43"
$ws.Range("E5").Value = "  +0.08%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5281"
$ws.Range("E6").Value = "  +1.76%  "

$ws.Range("E7").Value = "  -0.18%  "

$ws.Range("E8").Value = "  -1.73%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06319"
$ws.Range("E9").Value = "  +0.61%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "20.37"
$ws.Range("E10").Value = "  -2.19%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07792"
$ws.Range("E11").Value = "  +0.24%  "

$ws.Range("E12").Value = "  +0.66%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.645.73"
$ws.Range("E13").Value = "  -0.59%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.879.34"

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.5491"
$ws.Range("E15").Value = "  +0.26%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0₅8215"
$ws.Range("E16").Value = "  +1.02%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "65.52"
$ws.Range("E17").Value = "  +0.77%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "26.088.72"
$ws.Range("E18").Value = "  -0.49%  "

$ws.Range("E19").Value = "  -0.19%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.579"
$ws.Range("E20").Value = "  -0.84%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "190.66"

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.06"
$ws.Range("E22").Value = "  -0.18%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.027"
$ws.Range("E23").Value = "  +0.24%  "

$ws.Range("E24").Value = "  -0.20%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "143.00"
$ws.Range("E25").Value = "  +2.56%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1238"
$ws.Range("E26").Value = "  +1.09%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.234"
$ws.Range("E27").Value = "  -0.91%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "16.05"
$ws.Range("E28").Value = "  -0.71%  "

$ws.Range("E29").Value = "  -0.84%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.05832"
$ws.Range("E30").Value = "  -1.78%  "

$ws.Range("E31").Value = "  -0.09%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.549"
$ws.Range("E32").Value = "  -0.01%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.266"
$ws.Range("E33").Value = "  -0.50%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.586"
$ws.Range("E34").Value = "  +0.10%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.412"

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.9449"
$ws.Range("E36").Value = "  -1.79%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.776"
$ws.Range("E37").Value = "  +0.26%  "

$ws.Range("E38").Value = "  +1.04%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01608"
$ws.Range("E39").Value = "  +0.93%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.8481"
$ws.Range("E40").Value = "  -0.55%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "104.55"
$ws.Range("E41").Value = "  +3.73%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.735"
$ws.Range("E43").Value = "  -4.94%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.029.85"
$ws.Range("E44").Value = "  +1.82%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.794.76"
$ws.Range("E45").Value = "  -0.37%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "57.14"
$ws.Range("E46").Value = "  +1.09%  "

$ws.Range("E47").Value = "  -0.51%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.4328"
$ws.Range("E48").Value = "  +1.66%  "

$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "7.834"
$ws.Range("E49").Value = "  -2.35%  "

$ws.Range("B50").Value = "Cronos"
$ws.Range("C50").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.05141"
$ws.Range("E50").Value = "  -0.50%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.458"
$ws.Range("E51").Value = "  +0.13%  "
